# Update "想去人数" (F column) values on sheets "展览" and "全部类型"
# to reflect refreshed counts from the generated data source.

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetAll     = $wb.Worksheets.Item("全部类型")

# Row -> new F value for "展览" sheet
$exhibitUpdates = @{
    2  = 204
    4  = 5268
    8  = 594
    9  = 557
    13 = 4213
    15 = 179
    16 = 159
    18 = 3308
    23 = 194
    24 = 116
    25 = 35
    27 = 70
    28 = 299
    29 = 28
}

foreach ($row in $exhibitUpdates.Keys) {
    $sheetExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new F value for "全部类型" sheet
$allUpdates = @{
    2  = 204
    5  = 5268
    9  = 594
    10 = 557
    14 = 4213
    16 = 179
    17 = 159
    19 = 3308
    24 = 194
    25 = 116
    26 = 35
    28 = 70
    29 = 299
    30 = 28
}

foreach ($row in $allUpdates.Keys) {
    $sheetAll.Range("F$row").Value = $allUpdates[$row]
}
